$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 287.47223
$ws.Range("J17").Value = 287.47223
$ws.Range("L17").Value = 862.41669
$ws.Range("N17").Value = -1198.41669
$ws.Range("H98").Value = 882
$ws.Range("I98").Value = 900
$ws.Range("J98").Value = 756
$ws.Range("K98").Value = 900
$ws.Range("L98").Value = 756
$ws.Range("M98").Value = 598
$ws.Range("N98").Value = -3752
$ws.Range("H122").Value = 882
$ws.Range("I122").Value = 900
$ws.Range("J122").Value = 756
$ws.Range("K122").Value = 2700
$ws.Range("L122").Value = 2268
$ws.Range("M122").Value = -250
$ws.Range("N122").Value = -7168
$ws.Range("H127").Value = 2400.9607
$ws.Range("J127").Value = 2527.9375
$ws.Range("L127").Value = 7583.8125
$ws.Range("N127").Value = -17503.8125
$ws.Range("H132").Value = 4634005
$ws.Range("I132").Value = 5686806.5
$ws.Range("J132").Value = 1678.9
$ws.Range("K132").Value = 17060419.5
$ws.Range("L132").Value = 5036.700000000001
$ws.Range("M132").Value = -17057889.5
$ws.Range("N132").Value = -10096.7
$ws.Range("H137").Value = 979.1739
$ws.Range("J137").Value = 997.3333
$ws.Range("L137").Value = 2991.9999
$ws.Range("N137").Value = -8091.9999
$ws.Range("H138").Value = 4208.107
$ws.Range("I138").Value = 2028.4375
$ws.Range("J138").Value = 5079.975
$ws.Range("K138").Value = 6085.3125
$ws.Range("L138").Value = 15239.925
$ws.Range("M138").Value = -945.3125
$ws.Range("N138").Value = -25519.925

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2313.28
$ws.Range("I74").Value = 1266.15
$ws.Range("J74").Value = 6501.8
$ws.Range("K74").Value = 1266.15
$ws.Range("L74").Value = 6501.8
$ws.Range("M74").Value = -392.1500000000001
$ws.Range("N74").Value = -8249.799999999999
$ws.Range("H77").Value = 2313.28
$ws.Range("I77").Value = 1266.15
$ws.Range("J77").Value = 6501.8
$ws.Range("K77").Value = 6330.75
$ws.Range("L77").Value = 32509
$ws.Range("M77").Value = -1962.75
$ws.Range("N77").Value = -41245
$ws.Range("H102").Value = 64563.125
$ws.Range("I102").Value = 93007.27
$ws.Range("J102").Value = 1986
$ws.Range("K102").Value = 93007.27
$ws.Range("L102").Value = 1986
$ws.Range("M102").Value = -91385.27
$ws.Range("N102").Value = -5230
$ws.Range("H132").Value = 1498.2778
$ws.Range("I132").Value = 1185.8966
$ws.Range("K132").Value = 3557.6898
$ws.Range("M132").Value = -1027.6898

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 501747.75
$ws.Range("I105").Value = 502490
$ws.Range("J105").Value = 501005.5
$ws.Range("K105").Value = 502490
$ws.Range("L105").Value = 501005.5
$ws.Range("M105").Value = -500743
$ws.Range("N105").Value = -504499.5
$ws.Range("H107").Value = 30341832
$ws.Range("I107").Value = 47678020
$ws.Range("J107").Value = 3501.5
$ws.Range("K107").Value = 47678020
$ws.Range("L107").Value = 3501.5
$ws.Range("M107").Value = -47676100
$ws.Range("N107").Value = -7341.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17972.191
$ws.Range("I31").Value = 26959.744
$ws.Range("J31").Value = 3367.4167
$ws.Range("K31").Value = 26959.744
$ws.Range("L31").Value = 3367.4167
$ws.Range("M31").Value = -26664.744
$ws.Range("N31").Value = -3957.4167
$ws.Range("H34").Value = 17972.191
$ws.Range("I34").Value = 26959.744
$ws.Range("J34").Value = 3367.4167
$ws.Range("K34").Value = 26959.744
$ws.Range("L34").Value = 3367.4167
$ws.Range("M34").Value = -26757.744
$ws.Range("N34").Value = -3771.4167
$ws.Range("H107").Value = 871.1667
$ws.Range("I107").Value = 948.0909
$ws.Range("K107").Value = 948.0909
$ws.Range("M107").Value = 971.9091

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 104.21429
$ws.Range("J38").Value = 198.42857
$ws.Range("L38").Value = 595.28571
$ws.Range("N38").Value = -1289.28571
$ws.Range("H121").Value = 8151.263
$ws.Range("I121").Value = 5258.778
$ws.Range("J121").Value = 10754.5
$ws.Range("K121").Value = 15776.334
$ws.Range("L121").Value = 32263.5
$ws.Range("M121").Value = -14466.334
$ws.Range("N121").Value = -34883.5
$ws.Range("H122").Value = 8346.615
$ws.Range("I122").Value = 457
$ws.Range("J122").Value = 17551.166
$ws.Range("K122").Value = 4113
$ws.Range("L122").Value = 157960.494
$ws.Range("M122").Value = -1663
$ws.Range("N122").Value = -162860.494
$ws.Range("H131").Value = 799.37
$ws.Range("J131").Value = 820.4396
$ws.Range("L131").Value = 2461.3188
$ws.Range("N131").Value = -12541.3188

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1263060.1
$ws.Range("I107").Value = 450
$ws.Range("J107").Value = 3367410.2
$ws.Range("K107").Value = 450
$ws.Range("L107").Value = 3367410.2
$ws.Range("M107").Value = 1470
$ws.Range("N107").Value = -3371250.2
$ws.Range("H122").Value = 5236.7
$ws.Range("I122").Value = 4171
$ws.Range("K122").Value = 12513
$ws.Range("M122").Value = -10063

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1012596
$ws.Range("I46").Value = 325
$ws.Range("J46").Value = 2531002.5
$ws.Range("K46").Value = 325
$ws.Range("L46").Value = 2531002.5
$ws.Range("M46").Value = -137
$ws.Range("N46").Value = -2531378.5
$ws.Range("H68").Value = 3715.5833
$ws.Range("I68").Value = 1374.75
$ws.Range("J68").Value = 4886
$ws.Range("K68").Value = 1374.75
$ws.Range("L68").Value = 4886
$ws.Range("M68").Value = -625.75
$ws.Range("N68").Value = -6384
$ws.Range("H71").Value = 3715.5833
$ws.Range("I71").Value = 1374.75
$ws.Range("J71").Value = 4886
$ws.Range("K71").Value = 6873.75
$ws.Range("L71").Value = 24430
$ws.Range("M71").Value = -3129.75
$ws.Range("N71").Value = -31918
$ws.Range("H82").Value = 1566.5714
$ws.Range("I82").Value = 1499.8889
$ws.Range("K82").Value = 1499.8889
$ws.Range("M82").Value = -1138.8889
$ws.Range("H85").Value = 1566.5714
$ws.Range("I85").Value = 1499.8889
$ws.Range("K85").Value = 1499.8889
$ws.Range("M85").Value = -251.8888999999999
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H132").Value = 6219
$ws.Range("J132").Value = 3686.1428
$ws.Range("L132").Value = 11058.4284
$ws.Range("N132").Value = -16118.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7145257
$ws.Range("I62").Value = 50000000
$ws.Range("J62").Value = 2800
$ws.Range("K62").Value = 50000000
$ws.Range("L62").Value = 2800
$ws.Range("M62").Value = -49999376
$ws.Range("N62").Value = -4048
$ws.Range("H65").Value = 7145257
$ws.Range("I65").Value = 50000000
$ws.Range("J65").Value = 2800
$ws.Range("K65").Value = 250000000
$ws.Range("L65").Value = 14000
$ws.Range("M65").Value = -249996880
$ws.Range("N65").Value = -20240
$ws.Range("H132").Value = 4505.375
$ws.Range("I132").Value = 6448
$ws.Range("K132").Value = 19344
$ws.Range("M132").Value = -16814
$ws.Range("H136").Value = 1234.225
$ws.Range("I136").Value = 458.1111
$ws.Range("K136").Value = 1374.3333
$ws.Range("M136").Value = 1175.6667
